$wb = $excel.ActiveWorkbook

# Datetime (handback) stamps recorded for each locale's handback report
$handbackDateTimes = @{ "zh-cn" = "2016-03-02 06:35:37"; "de-de" = "2016-03-02 06:35:56" }

# The Overview sheet mirrors the same "Status" text for each file/locale pair, so it
# needs to be updated to the new handed-back status as well.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Capture the existing hyperlink info (address + display text) for columns A and C on
    # rows 2 and 3 before we touch anything, so we can duplicate them into new columns E/F.
    $links = @{}
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range().Address()
        $links[$addr] = @{ Address = $h.Address(); Display = $h.TextToDisplay() }
    }

    $a2 = $links["$" + "A$" + "2"]
    $c2 = $links["$" + "C$" + "2"]
    $a3 = $links["$" + "A$" + "3"]
    $c3 = $links["$" + "C$" + "3"]
    $a4 = $links["$" + "A$" + "4"]

    # Update the status text for the handed-off rows: they have now been handed back.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Record the handback datetime (previously the default epoch placeholder).
    $ws.Range("G2").Value = $handbackDateTimes[$sheetName]
    $ws.Range("G3").Value = $handbackDateTimes[$sheetName]

    # Rebuild every hyperlink on the sheet so the new Latest Target File (E) / Latest
    # Handback File (F) columns sit alongside the existing Source File Name (A) /
    # Latest Handoff File (C) columns, in natural left-to-right, top-to-bottom order.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $a2.Address, [Type]::Missing, [Type]::Missing, $a2.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $c2.Address, [Type]::Missing, [Type]::Missing, $c2.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), $a2.Address, [Type]::Missing, [Type]::Missing, $a2.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $c2.Address, [Type]::Missing, [Type]::Missing, $c2.Display) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $a3.Address, [Type]::Missing, [Type]::Missing, $a3.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $c3.Address, [Type]::Missing, [Type]::Missing, $c3.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E3"), $a3.Address, [Type]::Missing, [Type]::Missing, $a3.Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $c3.Address, [Type]::Missing, [Type]::Missing, $c3.Display) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $a4.Address, [Type]::Missing, [Type]::Missing, $a4.Display) | Out-Null
}

Write-Host "Handback report generated."
